$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 223 ("Romain Thunet" / "01:26:21" entry dated 2025-08-05) was removed
# from the sheet; every following row shifts up by one (224->223, ... 230->229).
$ws.Rows.Item(223).Delete()

# Mirror the author's final selection (the now-deleted row's former slot,
# i.e. the row that slid up into row 223, is left selected as a whole row).
$ws.Rows.Item(223).Select() | Out-Null
